# Update market/profit data values on several sheets, as produced by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 2331814.8
$ws.Range("I33").Value = 827.8125
$ws.Range("J33").Value = 12987755
$ws.Range("K33").Value = 827.8125
$ws.Range("L33").Value = 12987755
$ws.Range("M33").Value = -598.8125
$ws.Range("N33").Value = -12988213
# Row 40
$ws.Range("H40").Value = 11365535
$ws.Range("I40").Value = 1912.4138
$ws.Range("J40").Value = 33335206
$ws.Range("K40").Value = 1912.4138
$ws.Range("L40").Value = 33335206
$ws.Range("M40").Value = -1737.4138
$ws.Range("N40").Value = -33335556
# Row 112
$ws.Range("H112").Value = 21979384
$ws.Range("J112").Value = 27212414
$ws.Range("L112").Value = 81637242
$ws.Range("N112").Value = -81639458
# Row 137
$ws.Range("H137").Value = 1210.119
$ws.Range("I137").Value = 1145.7587
$ws.Range("J137").Value = 1353.6923
$ws.Range("K137").Value = 3437.2761
$ws.Range("L137").Value = 4061.0769
$ws.Range("M137").Value = -887.2761
$ws.Range("N137").Value = -9161.0769
# Row 138
$ws.Range("H138").Value = 3256.2415
$ws.Range("I138").Value = 1675.2258
$ws.Range("J138").Value = 4131.4463
$ws.Range("K138").Value = 5025.6774
$ws.Range("L138").Value = 12394.3389
$ws.Range("M138").Value = 114.3226000000004
$ws.Range("N138").Value = -22674.3389

$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 4633243
$ws.Range("I102").Value = 6174491
$ws.Range("J102").Value = 9500
$ws.Range("K102").Value = 6174491
$ws.Range("L102").Value = 9500
$ws.Range("M102").Value = -6172869
$ws.Range("N102").Value = -12744

$ws = $wb.Worksheets.Item("BSM")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
# Row 57
$ws.Range("H57").Value = 38000
$ws.Range("J57").Value = 38000
$ws.Range("L57").Value = 38000
$ws.Range("N57").Value = -39440
# Row 64
$ws.Range("H64").Value = 279
$ws.Range("I64").Value = 263
$ws.Range("J64").Value = 285.85715
$ws.Range("K64").Value = 263
$ws.Range("L64").Value = 285.85715
$ws.Range("M64").Value = -38
$ws.Range("N64").Value = -735.85715
# Row 67
$ws.Range("H67").Value = 279
$ws.Range("I67").Value = 263
$ws.Range("J67").Value = 285.85715
$ws.Range("K67").Value = 263
$ws.Range("L67").Value = 285.85715
$ws.Range("M67").Value = 517
$ws.Range("N67").Value = -1845.85715
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
# Row 133
$ws.Range("H133").Value = 44571.43
$ws.Range("J133").Value = 45333.332
$ws.Range("L133").Value = 45333.332
$ws.Range("N133").Value = -55453.332
# Row 134
$ws.Range("H134").Value = 50413.145
$ws.Range("I134").Value = 2982.2104
$ws.Range("J134").Value = 501007
$ws.Range("K134").Value = 8946.6312
$ws.Range("L134").Value = 1503021
$ws.Range("M134").Value = -6411.6312
$ws.Range("N134").Value = -1508091
# Row 135
$ws.Range("H135").Value = 70397.14
$ws.Range("J135").Value = 70397.14
$ws.Range("L135").Value = 70397.14
$ws.Range("N135").Value = -80537.14
# Row 136
$ws.Range("H136").Value = 38000
$ws.Range("J136").Value = 38000
$ws.Range("L136").Value = 38000
$ws.Range("N136").Value = -48200
# Row 137
$ws.Range("H137").Value = 59780
$ws.Range("J137").Value = 59780
$ws.Range("L137").Value = 59780
$ws.Range("N137").Value = -69980

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9809526
$ws.Range("I31").Value = 1787.3334
$ws.Range("J31").Value = 16674942
$ws.Range("K31").Value = 1787.3334
$ws.Range("L31").Value = 16674942
$ws.Range("M31").Value = -1492.3334
$ws.Range("N31").Value = -16675532
# Row 34
$ws.Range("H34").Value = 9809526
$ws.Range("I34").Value = 1787.3334
$ws.Range("J34").Value = 16674942
$ws.Range("K34").Value = 1787.3334
$ws.Range("L34").Value = 16674942
$ws.Range("M34").Value = -1585.3334
$ws.Range("N34").Value = -16675346
# Row 58
$ws.Range("H58").Value = 4310278
$ws.Range("I58").Value = 4903129
$ws.Range("J58").Value = 1430715.9
$ws.Range("K58").Value = 4903129
$ws.Range("L58").Value = 1430715.9
$ws.Range("M58").Value = -4902926
$ws.Range("N58").Value = -1431121.9
# Row 94
$ws.Range("H94").Value = 4508
$ws.Range("J94").Value = 4687.8
$ws.Range("L94").Value = 4687.8
$ws.Range("N94").Value = -5589.8
# Row 136
$ws.Range("H136").Value = 4310278
$ws.Range("I136").Value = 4903129
$ws.Range("J136").Value = 1430715.9
$ws.Range("K136").Value = 14709387
$ws.Range("L136").Value = 4292147.699999999
$ws.Range("M136").Value = -14706837
$ws.Range("N136").Value = -4297247.699999999
# Row 140
$ws.Range("H140").Value = 31038.46
$ws.Range("J140").Value = 31038.46
$ws.Range("L140").Value = 31038.46
$ws.Range("N140").Value = -41398.46

$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 4893
$ws.Range("I81").Value = 982.5
$ws.Range("J81").Value = 7500
$ws.Range("K81").Value = 2947.5
$ws.Range("L81").Value = 22500
$ws.Range("M81").Value = -1824.5
$ws.Range("N81").Value = -24746
# Row 84
$ws.Range("H84").Value = 4893
$ws.Range("I84").Value = 982.5
$ws.Range("J84").Value = 7500
$ws.Range("K84").Value = 8842.5
$ws.Range("L84").Value = 67500
$ws.Range("M84").Value = -3226.5
$ws.Range("N84").Value = -78732
# Row 113
$ws.Range("H113").Value = 2500500
$ws.Range("J113").Value = 476714.75
$ws.Range("L113").Value = 1430144.25
$ws.Range("N113").Value = -1434484.25
# Row 122
$ws.Range("H122").Value = 6651.6313
$ws.Range("I122").Value = 1107.2
$ws.Range("J122").Value = 12812.111
$ws.Range("K122").Value = 9964.800000000001
$ws.Range("L122").Value = 115308.999
$ws.Range("M122").Value = -7514.800000000001
$ws.Range("N122").Value = -120208.999
# Row 124
$ws.Range("H124").Value = 4999.8335
$ws.Range("I124").Value = 2499.75
$ws.Range("J124").Value = 10000
$ws.Range("K124").Value = 7499.25
$ws.Range("L124").Value = 30000
$ws.Range("M124").Value = -2589.25
$ws.Range("N124").Value = -39820

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 18521328
$ws.Range("I132").Value = 41668464
$ws.Range("J132").Value = 3621.8
$ws.Range("K132").Value = 125005392
$ws.Range("L132").Value = 10865.4
$ws.Range("M132").Value = -125002862
$ws.Range("N132").Value = -15925.4

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 11168183
$ws.Range("I81").Value = 1292.3846
$ws.Range("J81").Value = 40202100
$ws.Range("K81").Value = 2584.7692
$ws.Range("L81").Value = 80404200
$ws.Range("M81").Value = -1523.7692
$ws.Range("N81").Value = -80406322
# Row 84
$ws.Range("H84").Value = 11168183
$ws.Range("I84").Value = 1292.3846
$ws.Range("J84").Value = 40202100
$ws.Range("K84").Value = 12923.846
$ws.Range("L84").Value = 402021000
$ws.Range("M84").Value = -7619.846000000001
$ws.Range("N84").Value = -402031608
# Row 96
$ws.Range("H96").Value = 2821.2
$ws.Range("I96").Value = 2821.2
$ws.Range("K96").Value = 2821.2
$ws.Range("M96").Value = -1448.2
# Row 100
$ws.Range("H100").Value = 483.33334
$ws.Range("I100").Value = 475
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 950
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -409
$ws.Range("N100").Value = -2082
